$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valentin")

$ws.Range("A26").Value = "World of Warcraft Chronicle Volume 4"
$ws.Range("B26").Value = "https://m.media-amazon.com/images/I/816zRN+6uYL._SY466_.jpg"
$ws.Range("C26").Value = "https://www.amazon.de/-/en/gp/product/1506731910/ref=ox_sc_saved_image_1?smid=A3OJWAJQNSBARP&psc=1"
$ws.Range("D26").Value = "47 EUR"

$ws.Range("A27").Value = "Interstellar Blu-ray 4K HDR"
$ws.Range("B27").Value = "https://m.media-amazon.com/images/I/71k6YEIUXsL._SX342_.jpg"
$ws.Range("C27").Value = "https://www.amazon.de/-/en/gp/product/B075F3N9N5/ref=ox_sc_saved_image_4?smid=A3JWKAKR8XB7XF&psc=1"
$ws.Range("D27").Value = "18 EUR"

$ws.Range("A28").Value = "Jaws Blu-ray 4K"
$ws.Range("B28").Value = "https://m.media-amazon.com/images/I/81et3J4z9zL._SY445_.jpg"
$ws.Range("C28").Value = "https://www.amazon.de/-/en/gp/product/B0877NWJZF/ref=ox_sc_saved_image_7?smid=A30832IF5KZPY9&psc=1"
$ws.Range("D28").Value = "24 EUR"

$ws.Range("A29").Value = "Hario V60 Filter Holder"
$ws.Range("B29").Value = "https://m.media-amazon.com/images/I/61OJDLw1I5L._AC_SX679_.jpg"
$ws.Range("C29").Value = "https://www.amazon.de/-/en/gp/product/B00B7XJTP6/ref=ox_sc_saved_image_10?smid=ASHJXH5NF2K5S&psc=1"
$ws.Range("D29").Value = "28 EUR"

$ws.Range("A30").Value = "Blackwing Palomino 602 Crayons"
$ws.Range("B30").Value = "https://m.media-amazon.com/images/I/61LScqtldhL._AC_SX679_.jpg"
$ws.Range("C30").Value = "https://www.amazon.de/-/en/gp/product/B006YYPIUI/ref=ox_sc_saved_image_28?smid=A1JA9WU0P0W558&psc=1"
$ws.Range("D30").Value = "35 EUR"

$ws.Range("A31").Select() | Out-Null
